$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nbsp = [char]0x00A0
$batsman = "Pat Cummins" + $nbsp

# New match rows to append starting at row 13 (columns A..K):
# venue, date, result, ownTeam, oppTeam, batsman, totalRuns, totalBalls, total4s, total6s, sr
$rows = @(
    @(" Sharjah",      " October 26 2020",   "Kings XI won by 8 wickets (with 7 balls remaining)", "Kolkata Knight Riders", "Kings XI Punjab",             $batsman, "1",  "8",  "0", "0", "12.50"),
    @(" Abu Dhabi",     " September 23 2020", "Mumbai won by 49 runs",                               "Kolkata Knight Riders", "Mumbai Indians",              $batsman, "33", "12", "1", "4", "275.00"),
    @(" Abu Dhabi",     " October 24 2020",   "KKR won by 59 runs",                                  "Kolkata Knight Riders", "Delhi Capitals",              $batsman, "0",  "0",  "0", "0", "-"),
    @(" Dubai (DSC)",   " September 30 2020", "KKR won by 37 runs",                                  "Kolkata Knight Riders", "Rajasthan Royals",            $batsman, "12", "10", "1", "0", "120.00"),
    @(" Abu Dhabi",     " October 16 2020",   "Mumbai won by 8 wickets (with 19 balls remaining)",  "Kolkata Knight Riders", "Mumbai Indians",              $batsman, "53", "36", "5", "2", "147.22"),
    @(" Dubai (DSC)",   " November 01 2020",  "KKR won by 60 runs",                                  "Kolkata Knight Riders", "Rajasthan Royals",            $batsman, "15", "11", "0", "1", "136.36"),
    @(" Abu Dhabi",     " October 10 2020",   "KKR won by 2 runs",                                   "Kolkata Knight Riders", "Kings XI Punjab",             $batsman, "5",  "4",  "0", "0", "125.00"),
    @(" Sharjah",       " October 03 2020",   "Capitals won by 18 runs",                              "Kolkata Knight Riders", "Delhi Capitals",              $batsman, "5",  "4",  "1", "0", "125.00"),
    @(" Abu Dhabi",     " October 21 2020",   "RCB won by 8 wickets (with 39 balls remaining)",     "Kolkata Knight Riders", "Royal Challengers Bangalore", $batsman, "4",  "17", "0", "0", "23.52"),
    @(" Sharjah",       " October 12 2020",   "RCB won by 82 runs",                                   "Kolkata Knight Riders", "Royal Challengers Bangalore", $batsman, "1",  "3",  "0", "0", "33.33"),
    @(" Abu Dhabi",     " October 07 2020",   "KKR won by 10 runs",                                   "Kolkata Knight Riders", "Chennai Super Kings",         $batsman, "17", "9",  "1", "1", "188.88")
)

$startRow = 13
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    for ($col = 1; $col -le $data.Count; $col++) {
        $value = $data[$col - 1]
        $cell = $ws.Cells.Item($r, $col)
        if ($col -ge 7 -and $col -le 11 -and $value -match '^-?[0-9]+(\.[0-9]+)?$') {
            # Numeric-looking values (totalRuns, totalBalls, total4s, total6s, sr) are
            # stored as text in this sheet, so force text using a leading apostrophe.
            $cell.Value = "'" + $value
        } else {
            $cell.Value = $value
        }
    }
}
